$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.567.67"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").Value = "3.477.76"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.04"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.51"
$ws.Range("E6").Value = "  +2.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D8").Value = "3.477.17"
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.585"
$ws.Range("E9").Value = "  +5.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.28"
$ws.Range("E10").Value = "  -4.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.125"
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.442"
$ws.Range("E12").Value = "  -1.54%  "
$ws.Range("D13").Value = "4.071.29"
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("E14").Value = "  -1.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000194"
$ws.Range("E15").Value = "  -2.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.55"
$ws.Range("E16").Value = "  +2.13%  "
$ws.Range("D17").Value = "65.523.03"
$ws.Range("E17").Value = "  +1.08%  "
$ws.Range("D18").Value = "3.471.07"
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.40"
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.23"
$ws.Range("E20").Value = "  -1.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "390.08"
$ws.Range("E21").Value = "  -2.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.23"
$ws.Range("E22").Value = "  -3.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.548"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.56"
$ws.Range("E24").Value = "  +0.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.997"
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000123"
$ws.Range("E26").Value = "  +0.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.56"
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.179"
$ws.Range("E28").Value = "  -1.07%  "
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("E30").Value = "  +6.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.42"
$ws.Range("E31").Value = "  +3.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.05"
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.63"
$ws.Range("E33").Value = "  -1.06%  "
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.49"
$ws.Range("E34").Value = "  -4.08%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.09"
$ws.Range("E36").Value = "  +0.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.53"
$ws.Range("E37").Value = "  +3.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "162.25"
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.95"
$ws.Range("E39").Value = "  +4.27%  "
$ws.Range("D40").Value = "3.047.65"
$ws.Range("E40").Value = "  +4.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0770"
$ws.Range("E41").Value = "  -1.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.94"
$ws.Range("E42").Value = "  -2.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0320"
$ws.Range("E43").Value = "  -1.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.50"
$ws.Range("E44").Value = "  +1.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.65"
$ws.Range("E45").Value = "  +1.65%  "
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.77"
$ws.Range("E47").Value = "  +8.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.11"
$ws.Range("E48").Value = "  +1.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.20"
$ws.Range("E49").Value = "  +0.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.68"
$ws.Range("E50").Value = "  +1.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "309.65"
$ws.Range("E51").Value = "  +3.26%  "
